$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-0.8860551090425615, 1.2466813679246638, 1.128640659030754, 3.1296047265008546)"
$ws.Range("C2").Value = "NCT(3.0649655872736417, 2.630029135862989, -3.5938826991225143, 4.220272025460711)"
$ws.Range("D2").Value = "EXN(2.944687840166079, 0.3969455356773065, 1.5781073873341749)"
$ws.Range("E2").Value = "NCT(2.166288802748673, 2.225071355342058, -0.008614583288274962, 2.993485620643642)"
